$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number; force text so Excel doesn't coerce it to a number
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 25.05.2024"

$ws.Range("B6").Value = "28.05."
$ws.Range("C6").Value = "29.05."
$ws.Range("D6").Value = "AMAZON.DE MKTPLC EU BTAUXH"
$ws.Range("E6").Value = "195,89-"

$ws.Range("B7").Value = "01.06."
$ws.Range("C7").Value = "02.06."
$ws.Range("D7").Value = "KARTENZ./01.06 ALDI SUED RO"
$ws.Range("E7").Value = "89,34-"

$ws.Range("B8").Value = "04.06."
$ws.Range("C8").Value = "05.06."
$ws.Range("D8").Value = "KARTENZ./04.06 LIDL RO"
$ws.Range("E8").Value = "84,31-"

$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 49496749"
$ws.Range("E9").Value = "37,60-"

$ws.Range("B10").Value = "07.06."
$ws.Range("C10").Value = "08.06."
$ws.Range("D10").Value = "BURGER KING Soltau"
$ws.Range("E10").Value = "21,93-"

$ws.Range("D12").Value = "KONTOSTAND AM 12.06.2024"
$ws.Range("E12").Value = "429,07-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 20.06.2024"
